$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new data row (row 4) mirroring the existing rows' layout.
$ws.Cells.Item(4, 1).Value = 42606.881168981483
$ws.Cells.Item(4, 2).Value = 30
$ws.Cells.Item(4, 3).Value = 65
$ws.Cells.Item(4, 4).Value = 34
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 4753
$ws.Cells.Item(4, 8).Value = 2294
$ws.Cells.Item(4, 9).Value = 128
$ws.Cells.Item(4, 10).Value = 21
$ws.Cells.Item(4, 11).Value = 11
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).Value = 0
$ws.Cells.Item(4, 14).Value = "Named"

# Match the date-style formatting used by the other rows in column A
# (same builtin format as the existing date cells, so it reuses the same style).
$ws.Cells.Item(4, 1).NumberFormat = "m/d/yy h:mm"
